$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86: 'Filling in the Blanks' / 'Enchanted Aurum Regis Ink'
$ws.Range("H86").Value = 726193.1
$ws.Range("I86").Value = 3800
$ws.Range("J86").Value = 1689384
$ws.Range("K86").Value = 3800
$ws.Range("L86").Value = 1689384
$ws.Range("M86").Value = -2677
$ws.Range("N86").Value = -1691630

# Row 89: 'Ink into Antiquity (L)' / 'Enchanted Aurum Regis Ink'
$ws.Range("H89").Value = 726193.1
$ws.Range("I89").Value = 3800
$ws.Range("J89").Value = 1689384
$ws.Range("K89").Value = 19000
$ws.Range("L89").Value = 8446920
$ws.Range("M89").Value = -13384
$ws.Range("N89").Value = -8458152

$ws = $wb.Worksheets.Item("ARM")
# Row 32: 'Ingot We Trust' / 'Steel Ingot'
$ws.Range("H32").Value = 1629119
$ws.Range("I32").Value = 1989941.5
$ws.Range("J32").Value = 5417.4287
$ws.Range("K32").Value = 1989941.5
$ws.Range("L32").Value = 5417.4287
$ws.Range("M32").Value = -1989654.5
$ws.Range("N32").Value = -5991.4287

# Row 74: 'As the Bolt Flies' / 'Titanium Nugget'
$ws.Range("H74").Value = 26054730
$ws.Range("I74").Value = 23999322
$ws.Range("J74").Value = 31374614
$ws.Range("K74").Value = 23999322
$ws.Range("L74").Value = 31374614
$ws.Range("M74").Value = -23998448
$ws.Range("N74").Value = -31376362

# Row 77: 'Heavy Metal Banned (L)' / 'Titanium Nugget'
$ws.Range("H77").Value = 26054730
$ws.Range("I77").Value = 23999322
$ws.Range("J77").Value = 31374614
$ws.Range("K77").Value = 119996610
$ws.Range("L77").Value = 156873070
$ws.Range("M77").Value = -119992242
$ws.Range("N77").Value = -156881806

# Row 132: "Don't Bore Me, Ore Me" / 'Mountain Chromite Ingot'
$ws.Range("H132").Value = 16670486
$ws.Range("I132").Value = 13378919
$ws.Range("J132").Value = 25557716
$ws.Range("K132").Value = 40136757
$ws.Range("L132").Value = 76673148
$ws.Range("M132").Value = -40134227
$ws.Range("N132").Value = -76678208

$ws = $wb.Worksheets.Item("BSM")
# Row 86: 'Through Thick and Thin' / 'Adamantite Nugget'
$ws.Range("H86").Value = 1936.91
$ws.Range("I86").Value = 1944.8062
$ws.Range("J86").Value = 1550
$ws.Range("K86").Value = 1944.8062
$ws.Range("L86").Value = 1550
$ws.Range("M86").Value = -821.8062
$ws.Range("N86").Value = -3796

# Row 89: 'Piercing Eyes Deserve Piercing Shafts (L)' / 'Adamantite Nugget'
$ws.Range("H89").Value = 1936.91
$ws.Range("I89").Value = 1944.8062
$ws.Range("J89").Value = 1550
$ws.Range("K89").Value = 9724.030999999999
$ws.Range("L89").Value = 7750
$ws.Range("M89").Value = -4108.030999999999
$ws.Range("N89").Value = -18982

# Row 94: 'High Steal' / 'High Steel Nugget'
$ws.Range("H94").Value = 1712.52
$ws.Range("I94").Value = 1223.5
$ws.Range("J94").Value = 2970
$ws.Range("K94").Value = 1223.5
$ws.Range("L94").Value = 2970
$ws.Range("M94").Value = -772.5
$ws.Range("N94").Value = -3872

# Row 134: 'Ruthenium Supremium' / 'Ruthenium Ingot'
$ws.Range("H134").Value = 33089464
$ws.Range("I134").Value = 38462844
$ws.Range("J134").Value = 9804821
$ws.Range("K134").Value = 115388532
$ws.Range("L134").Value = 29414463
$ws.Range("M134").Value = -115385997
$ws.Range("N134").Value = -29419533

$ws = $wb.Worksheets.Item("CRP")
# Row 31: 'Wall Not Found' / 'Walnut Lumber'
$ws.Range("H31").Value = 14086454
$ws.Range("I31").Value = 33334656
$ws.Range("J31").Value = 2404.2195
$ws.Range("K31").Value = 33334656
$ws.Range("L31").Value = 2404.2195
$ws.Range("M31").Value = -33334361
$ws.Range("N31").Value = -2994.2195

# Row 34: 'Armoires of the Rich and Famous' / 'Walnut Lumber'
$ws.Range("H34").Value = 14086454
$ws.Range("I34").Value = 33334656
$ws.Range("J34").Value = 2404.2195
$ws.Range("K34").Value = 33334656
$ws.Range("L34").Value = 2404.2195
$ws.Range("M34").Value = -33334454
$ws.Range("N34").Value = -2808.2195

# Row 62: 'Splinter in the Sewers' / 'Cedar Lumber'
$ws.Range("H62").Value = 3288.889
$ws.Range("I62").Value = 2450
$ws.Range("K62").Value = 2450
$ws.Range("M62").Value = -1826

# Row 65: 'The Lumber of Their Discontent (L)' / 'Cedar Lumber'
$ws.Range("H65").Value = 3288.889
$ws.Range("I65").Value = 2450
$ws.Range("K65").Value = 12250
$ws.Range("M65").Value = -9130

# Row 95: 'Standing on Ceremony' / 'High Steel Fork'
$ws.Range("H95").Value = 12238.857
$ws.Range("J95").Value = 12238.857
$ws.Range("L95").Value = 12238.857
$ws.Range("N95").Value = -17730.857

$ws = $wb.Worksheets.Item("CUL")
# Row 5: 'What a Sap' / 'Maple Syrup'
$ws.Range("H5").Value = 1795419.9
$ws.Range("I5").Value = 1374044.8
$ws.Range("J5").Value = 2778628.5
$ws.Range("K5").Value = 4122134.4
$ws.Range("L5").Value = 8335885.5
$ws.Range("M5").Value = -4122022.4
$ws.Range("N5").Value = -8336109.5

# Row 131: 'The Mountain Steeped' / 'Tsai tou Vounou'
$ws.Range("H131").Value = 134597.33
$ws.Range("J131").Value = 112916.664
$ws.Range("L131").Value = 338749.992
$ws.Range("N131").Value = -348829.992

# Row 132: 'More Mezcal' / 'Cooking Mezcal'
$ws.Range("H132").Value = 1238.4634
$ws.Range("I132").Value = 1021.64703
$ws.Range("J132").Value = 1392.0416
$ws.Range("K132").Value = 9194.823269999999
$ws.Range("L132").Value = 12528.3744
$ws.Range("M132").Value = -6664.823269999999
$ws.Range("N132").Value = -17588.3744

# Row 133: 'Friends Are Food' / 'Boiled Alpaca Steak'
$ws.Range("H133").Value = 2555
$ws.Range("I133").Value = 2352.7273
$ws.Range("J133").Value = 3000
$ws.Range("K133").Value = 7058.1819
$ws.Range("L133").Value = 9000
$ws.Range("M133").Value = -1998.1819
$ws.Range("N133").Value = -19120

# Row 135: 'Not-so-secret Ingredient' / 'Royal Maple Syrup'
$ws.Range("H135").Value = 1795419.9
$ws.Range("I135").Value = 1374044.8
$ws.Range("J135").Value = 2778628.5
$ws.Range("K135").Value = 12366403.2
$ws.Range("L135").Value = 25007656.5
$ws.Range("M135").Value = -12363868.2
$ws.Range("N135").Value = -25012726.5

# Row 138: 'Bring Me Your Tacos' / 'Tacos Al Pastor'
$ws.Range("H138").Value = 1919.1
$ws.Range("I138").Value = 662.5714
$ws.Range("J138").Value = 4851
$ws.Range("K138").Value = 1987.7142
$ws.Range("L138").Value = 14553
$ws.Range("M138").Value = 3152.2858
$ws.Range("N138").Value = -24833

# Row 139: 'Najoothie' / 'Wild Banana Blend'
$ws.Range("H139").Value = 75006.57000000001
$ws.Range("I139").Value = 250015
$ws.Range("J139").Value = 5003.2
$ws.Range("K139").Value = 750045
$ws.Range("L139").Value = 15009.6
$ws.Range("M139").Value = -744905
$ws.Range("N139").Value = -25289.6

# Row 140: 'Sweet, Sweet Bean Juice' / 'Mesquite Juice'
$ws.Range("H140").Value = 2265.3157
$ws.Range("I140").Value = 2389.5
$ws.Range("J140").Value = 2052.4285
$ws.Range("K140").Value = 7168.5
$ws.Range("L140").Value = 6157.2855
$ws.Range("M140").Value = -1988.5
$ws.Range("N140").Value = -16517.2855

$ws = $wb.Worksheets.Item("GSM")
# Row 52: "It's My Business to Know Things" / 'Red Coral Armillae'
$ws.Range("H52").Value = 12681.818
$ws.Range("J52").Value = 12681.818
$ws.Range("L52").Value = 12681.818
$ws.Range("N52").Value = -13199.818

$ws = $wb.Worksheets.Item("LTW")
# Row 7: 'Tan Before the Ban' / 'Leather'
$ws.Range("H7").Value = 2063.125
$ws.Range("I7").Value = 1928.5714
$ws.Range("J7").Value = 3005
$ws.Range("K7").Value = 1928.5714
$ws.Range("L7").Value = 3005
$ws.Range("M7").Value = -1816.5714
$ws.Range("N7").Value = -3229

# Row 40: 'Best Served Toad' / 'Toad Leather'
$ws.Range("H40").Value = 18523134
$ws.Range("I40").Value = 27781452
$ws.Range("J40").Value = 6499
$ws.Range("K40").Value = 27781452
$ws.Range("L40").Value = 6499
$ws.Range("M40").Value = -27781316
$ws.Range("N40").Value = -6771

# Row 126: 'Battered Books' / 'Saiga Leather'
$ws.Range("H126").Value = 2063.125
$ws.Range("I126").Value = 1928.5714
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 5785.7142
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -3315.7142
$ws.Range("N126").Value = -13955

# Row 127: 'Loyal Turncoat' / 'Saigaskin Coat of Fending'
$ws.Range("H127").Value = 54780.125
$ws.Range("J127").Value = 54780.125
$ws.Range("L127").Value = 54780.125
$ws.Range("N127").Value = -64700.125
